$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-16) of this sheet get permuted: each destination row receives
# the full set of varying values (columns D, I, J, K, L, M, N, O, P, Q) that
# previously belonged to a different source row. Columns A, B, C, E, F, G, H, R
# are identical across all rows and are left untouched.
#
# destination row -> source row (both refer to the original/"before" layout)
$map = @{
    2  = 7
    3  = 8
    4  = 3
    5  = 4
    6  = 15
    7  = 11
    8  = 2
    9  = 16
    10 = 6
    11 = 12
    12 = 9
    13 = 14
    14 = 13
    15 = 10
    16 = 5
}

# Columns whose values differ row to row and must be carried over by the permutation.
$cols = @("D", "I", "J", "K", "L", "M", "N", "O", "P", "Q")

# First snapshot all current ("before") values for every row so that writes to
# earlier rows don't clobber data still needed for later rows.
# NOTE: ".Value2" (not ".Value") is used to read, since the interop's ".Value"
# getter does not reliably return cell contents in this runtime.
$snapshot = @{}
for ($r = 2; $r -le 16; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write each destination row using the snapshot of its mapped source row.
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $srcVals[$col]
    }
}
